# Automatic update of files.
#
# 1) Column C ("Förändrad") gets bumped from serial date 45184 to 45186
#    for every data row (rows 2-79).
# 2) For the first 8 data rows (rows 2-9), the HYPERLINK() formulas in
#    columns S, T, V, W, X and Y gain a second "friendly name" argument
#    equal to the case/report code held in column A of that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Bump the "Förändrad" date column for every data row -----------
$ws.Range("C2:C79").Value = 45186

# --- 2) Add the friendly-name argument to the HYPERLINK formulas ------
$linkCols = @("S", "T", "V", "W", "X", "Y")

for ($row = 2; $row -le 9; $row++) {
    $code = $ws.Range("A" + $row).Text

    foreach ($col in $linkCols) {
        $cell = $ws.Range($col + $row)
        $formula = $cell.Formula
        if ($formula -and $formula -match '^=HYPERLINK\((.*)\)$') {
            $args = $Matches[1]
            $cell.Formula = '=HYPERLINK(' + $args + ', "' + $code + '")'
        }
    }
}
